$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "[1, 0, 1]"

$ws.Range("B3").Value = "[1, 0, 1]"

$ws.Range("C4").Value = "[0, 0, 1]"

$ws.Range("B5").Value = "[0, 0, 1]"
$ws.Range("C5").Value = "[1, 0, 0]"

$ws.Range("B6").Value = "[0, 1, 0]"
$ws.Range("C6").Value = "[1, 0, 1]"

$ws.Range("B7").Value = "[1, 1, 0]"
$ws.Range("C7").Value = "[0, 0, 1]"

$ws.Range("B8").Value = "[0, 1, 1]"
$ws.Range("C8").Value = "[1, 0, 0]"

$ws.Range("B9").Value = "[1, 1, 1]"
$ws.Range("C9").Value = "[0, 0, 0]"
